$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("manager_ids")

# Append new rows for the 2023/24 season
$ws.Range("A29").Value = 2564951
$ws.Range("B29").Value = "Andy"
$ws.Range("C29").Value = "Season 23-24"
$ws.Range("D29").Value = 12345

$ws.Range("A30").Value = 2565192
$ws.Range("B30").Value = "Spartak Kent"
$ws.Range("C30").Value = "Season 23-24"
$ws.Range("D30").Value = 12345
